$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").Value = ""
$ws.Range("H31").Value = 3268.8572
$ws.Range("I31").Value = 1980.5
$ws.Range("K31").Value = 5941.5
$ws.Range("M31").Value = -5711.5
$ws.Range("H62").Value = 1333.3334
$ws.Range("J62").Value = 2000
$ws.Range("L62").Value = 2000
$ws.Range("N62").Value = -3248
$ws.Range("H65").Value = 1333.3334
$ws.Range("J65").Value = 2000
$ws.Range("L65").Value = 10000
$ws.Range("N65").Value = -16240
$ws.Range("H86").Value = 9867.166999999999
$ws.Range("I86").Value = 14916.333
$ws.Range("J86").Value = 4818
$ws.Range("K86").Value = 14916.333
$ws.Range("L86").Value = 4818
$ws.Range("M86").Value = -13793.333
$ws.Range("N86").Value = -7064
$ws.Range("H89").Value = 9867.166999999999
$ws.Range("I89").Value = 14916.333
$ws.Range("J89").Value = 4818
$ws.Range("K89").Value = 74581.66500000001
$ws.Range("L89").Value = 24090
$ws.Range("M89").Value = -68965.66500000001
$ws.Range("N89").Value = -35322
$ws.Range("H103").Value = 781.8946999999999
$ws.Range("I103").Value = 1043.8889
$ws.Range("J103").Value = 546.1
$ws.Range("K103").Value = 3131.6667
$ws.Range("L103").Value = 1638.3
$ws.Range("M103").Value = -2545.6667
$ws.Range("N103").Value = -2810.3
$ws.Range("H138").Value = 2465.8718
$ws.Range("I138").Value = 2035.2222
$ws.Range("J138").Value = 2595.0667
$ws.Range("K138").Value = 6105.6666
$ws.Range("L138").Value = 7785.2001
$ws.Range("M138").Value = -965.6665999999996
$ws.Range("N138").Value = -18065.2001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1052465
$ws.Range("J2").Value = 1848.3334
$ws.Range("L2").Value = 1848.3334
$ws.Range("N2").Value = -2074.3334
$ws.Range("H45").Value = 4204.6665
$ws.Range("I45").Value = 5058.8237
$ws.Range("J45").Value = 2130.2856
$ws.Range("K45").Value = 5058.8237
$ws.Range("L45").Value = 2130.2856
$ws.Range("M45").Value = -4681.8237
$ws.Range("N45").Value = -2884.2856
$ws.Range("H88").Value = 2653.2856
$ws.Range("I88").Value = 2330.2
$ws.Range("J88").Value = 2832.7778
$ws.Range("K88").Value = 2330.2
$ws.Range("L88").Value = 2832.7778
$ws.Range("M88").Value = -1924.2
$ws.Range("N88").Value = -3644.7778
$ws.Range("H91").Value = 2653.2856
$ws.Range("I91").Value = 2330.2
$ws.Range("J91").Value = 2832.7778
$ws.Range("K91").Value = 2330.2
$ws.Range("L91").Value = 2832.7778
$ws.Range("M91").Value = -926.1999999999998
$ws.Range("N91").Value = -5640.7778
$ws.Range("H102").Value = 9092416
$ws.Range("I102").Value = 11112411
$ws.Range("K102").Value = 11112411
$ws.Range("M102").Value = -11110789
$ws.Range("H116").Value = 1052465
$ws.Range("J116").Value = 1848.3334
$ws.Range("L116").Value = 1848.3334
$ws.Range("N116").Value = -6436.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1052465
$ws.Range("J3").Value = 1848.3334
$ws.Range("L3").Value = 1848.3334
$ws.Range("N3").Value = -2076.3334
$ws.Range("H11").Value = 308
$ws.Range("I11").Value = 25
$ws.Range("K11").Value = 25
$ws.Range("M11").Value = 115
$ws.Range("H22").Value = 1140.5555
$ws.Range("I22").Value = 256.66666
$ws.Range("J22").Value = 1582.5
$ws.Range("K22").Value = 256.66666
$ws.Range("L22").Value = 1582.5
$ws.Range("M22").Value = -83.66665999999998
$ws.Range("N22").Value = -1928.5
$ws.Range("H54").Value = 44999
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").Value = ""
$ws.Range("H82").Value = 19902.143
$ws.Range("I82").Value = 7863
$ws.Range("J82").Value = 50000
$ws.Range("K82").Value = 7863
$ws.Range("L82").Value = 50000
$ws.Range("M82").Value = -7480
$ws.Range("N82").Value = -50766
$ws.Range("H85").Value = 19902.143
$ws.Range("I85").Value = 7863
$ws.Range("J85").Value = 50000
$ws.Range("K85").Value = 7863
$ws.Range("L85").Value = 50000
$ws.Range("M85").Value = -6537
$ws.Range("N85").Value = -52652
$ws.Range("H105").Value = 2658.842
$ws.Range("I105").Value = 2205.5
$ws.Range("K105").Value = 2205.5
$ws.Range("M105").Value = -458.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 4000026.5
$ws.Range("I4").Value = 33
$ws.Range("K4").Value = 33
$ws.Range("M4").Value = 79
$ws.Range("H86").Value = 10615.682
$ws.Range("I86").Value = 7505.125
$ws.Range("J86").Value = 12393.143
$ws.Range("K86").Value = 7505.125
$ws.Range("L86").Value = 12393.143
$ws.Range("M86").Value = -6382.125
$ws.Range("N86").Value = -14639.143
$ws.Range("H89").Value = 10615.682
$ws.Range("I89").Value = 7505.125
$ws.Range("J89").Value = 12393.143
$ws.Range("K89").Value = 37525.625
$ws.Range("L89").Value = 61965.715
$ws.Range("M89").Value = -31909.625
$ws.Range("N89").Value = -73197.715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 8812.25
$ws.Range("I26").Value = 260.42856
$ws.Range("J26").Value = 20784.8
$ws.Range("K26").Value = 781.28568
$ws.Range("L26").Value = 62354.39999999999
$ws.Range("M26").Value = -493.28568
$ws.Range("N26").Value = -62930.39999999999
$ws.Range("H59").Value = 3999
$ws.Range("I59").Value = 3999
$ws.Range("K59").Value = 11997
$ws.Range("M59").Value = -11457
$ws.Range("H113").Value = 500499.5
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = ""
$ws.Range("H122").Value = 1098.8
$ws.Range("I122").Value = 1124.25
$ws.Range("J122").Value = 997
$ws.Range("K122").Value = 10118.25
$ws.Range("L122").Value = 8973
$ws.Range("M122").Value = -7668.25
$ws.Range("N122").Value = -13873

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 83
$ws.Range("I2").Value = 38.285713
$ws.Range("J2").Value = 135.16667
$ws.Range("K2").Value = 38.285713
$ws.Range("L2").Value = 135.16667
$ws.Range("M2").Value = 74.714287
$ws.Range("N2").Value = -361.16667
$ws.Range("H22").Value = 100000
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = ""
$ws.Range("H58").Value = 15000
$ws.Range("J58").Value = 15000
$ws.Range("L58").Value = 15000
$ws.Range("N58").Value = -15554
$ws.Range("H126").Value = 3725.6365
$ws.Range("I126").Value = 3725.6365
$ws.Range("K126").Value = 11176.9095
$ws.Range("M126").Value = -8706.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = ""
$ws.Range("H7").Value = 2990.1428
$ws.Range("J7").Value = 3609.6
$ws.Range("L7").Value = 3609.6
$ws.Range("N7").Value = -3833.6
$ws.Range("H58").Value = 13574.625
$ws.Range("I58").Value = 10219.8
$ws.Range("J58").Value = 19166
$ws.Range("K58").Value = 10219.8
$ws.Range("L58").Value = 19166
$ws.Range("M58").Value = -9959.799999999999
$ws.Range("N58").Value = -19686
$ws.Range("H126").Value = 2990.1428
$ws.Range("J126").Value = 3609.6
$ws.Range("L126").Value = 10828.8
$ws.Range("N126").Value = -15768.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2099.8333
$ws.Range("I100").Value = 2263.4546
$ws.Range("J100").Value = 300
$ws.Range("K100").Value = 4526.9092
$ws.Range("L100").Value = 600
$ws.Range("M100").Value = -3985.9092
$ws.Range("N100").Value = -1682
